$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "starred" plants column (D) for several users.
# andrea: "Spiky boi" -> "Spiky boi, Tulips"
$ws.Range("D2").Value = "Spiky boi, Tulips"

# catherine: new starred plant "Peace lily"
$ws.Range("D5").Value = "Peace lily"

# mark: new starred plants "Peace lily, Ficus"
$ws.Range("D3").Value = "Peace lily, Ficus"

# ben: new starred plant "Aloe Vera"
$ws.Range("D8").Value = "Aloe Vera"

# Move the active selection to D5, matching the saved workbook view state.
$ws.Range("D5").Select()
